$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Bad Drivers) updates
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 118
$ws.Range("D3").Value = 98.90000000000001

# Row 4 (Totals) updates
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 118

# Row 12 - clear the Driver Vintage date value
$ws.Range("E12").ClearContents()

# Row 14 - update total samples
$ws.Range("B14").Value = 265400
